# Add new "intervention_type" column (K) with a value for every clinical
# trial row, mirroring the existing header/data styling conventions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1: copy the style of the existing header cell (A1) so the
# new header reuses the same bold / bordered / centered format, then set
# its text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "intervention_type"

# Data rows 2-49: plain values, no special styling (matches column J, etc.)
$interventionTypes = @{
    2 = "OTHER"
    3 = "OTHER"
    4 = "OTHER"
    5 = "OTHER"
    6 = "OTHER"
    7 = "DIAGNOSTIC_TEST"
    8 = "PROCEDURE"
    9 = "OTHER"
    10 = "PROCEDURE"
    11 = "DEVICE"
    12 = "OTHER"
    13 = "DIAGNOSTIC_TEST"
    14 = "DEVICE"
    15 = "PROCEDURE"
    16 = "DRUG"
    17 = "DRUG"
    18 = "DIAGNOSTIC_TEST"
    19 = "DIAGNOSTIC_TEST"
    20 = "DIAGNOSTIC_TEST"
    21 = "DEVICE"
    22 = "OTHER"
    23 = "OTHER"
    24 = "RADIATION"
    25 = "OTHER"
    26 = "DEVICE"
    27 = "OTHER"
    28 = "DEVICE"
    29 = "OTHER"
    30 = "OTHER"
    31 = "BIOLOGICAL"
    32 = "BEHAVIORAL"
    33 = "OTHER"
    34 = "OTHER"
    35 = "DEVICE"
    36 = "PROCEDURE"
    37 = "DIAGNOSTIC_TEST"
    38 = "DRUG"
    39 = "PROCEDURE"
    40 = "OTHER"
    41 = "OTHER"
    42 = "BIOLOGICAL"
    43 = "OTHER"
    44 = "DIAGNOSTIC_TEST"
    45 = "DIAGNOSTIC_TEST"
    46 = "DEVICE"
    47 = "PROCEDURE"
    48 = "DEVICE"
    49 = "PROCEDURE"
}

foreach ($row in $interventionTypes.Keys | Sort-Object { [int]$_ }) {
    $ws.Cells.Item([int]$row, 11).Value = $interventionTypes[$row]
}
